$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-15 20:28:06", 0.001),
    @("2023-12-15 20:29:34", 0.005600000000000002),
    @("2023-12-15 20:29:54", 0.0012),
    @("2023-12-15 20:30:01", 0.0004)
)

$startRow = 379
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
